# Add a "Myers Briggs" table to Sheet1: new header columns (Mind, Energy,
# Nature, Tactics, Identity) plus per-person data, wrapped in an Excel Table
# (ListObject) named "Table1" over A1:G7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Names"
# B1 ("Myers Brigs") already exists from the original sheet.
$ws.Range("C1").Value = "Mind"
$ws.Range("D1").Value = "Energy "
$ws.Range("E1").Value = "Nature"
$ws.Range("F1").Value = "Tactics"
$ws.Range("G1").Value = "Identity"

# --- Row 2 (Jae Sung Oh) keeps only the original two columns --------------
# (A2/B2 already populated - no Myers Briggs breakdown for this row.)

# --- Row 3 (Luke Green) ----------------------------------------------------
$ws.Range("C3").Value = "65% Introverted"
$ws.Range("D3").Value = "63% Observant"
$ws.Range("E3").Value = "56% Feeling"
$ws.Range("F3").Value = "75% Judging"
$ws.Range("G3").Value = "76% Assertive"

# --- Row 4 (Natalie Yelland-Hall) ------------------------------------------
$ws.Range("C4").Value = "100% Introverted"
$ws.Range("D4").Value = "63% Observant"
$ws.Range("E4").Value = "60% Feeling"
$ws.Range("F4").Value = "71% Judging"
$ws.Range("G4").Value = "65% Turbulent"
$ws.Range("A4").Font.Bold = $true

# --- Row 5 (Michael Seach) --------------------------------------------------
$ws.Range("C5").Value = "57% Extroverted"
$ws.Range("D5").Value = "68% Observant"
$ws.Range("E5").Value = "60% Feeling"
$ws.Range("F5").Value = "56% Prospecting"
$ws.Range("G5").Value = "56%Assertive"

# --- Row 6 (Bobbie Cole) -----------------------------------------------------
$ws.Range("C6").Value = "88% Extraverted"
$ws.Range("D6").Value = "64% Intuituve"
$ws.Range("E6").Value = "51% Feeling"
$ws.Range("F6").Value = "51% Prospecting"
$ws.Range("G6").Value = "53% Assertive"

# --- Row 7 (Phuong Quach) ---------------------------------------------------
$ws.Range("C7").Value = "79% Introverted"
$ws.Range("D7").Value = "56% Observant"
$ws.Range("E7").Value = "68% Thinking"
$ws.Range("F7").Value = "64% Prospecting"
$ws.Range("G7").Value = "97% Turbulent"

# --- Column widths for the new columns -------------------------------------
$ws.Range("C1").ColumnWidth = 16.0
$ws.Range("D1").ColumnWidth = 16.0
$ws.Range("E1").ColumnWidth = 18.166666666666668
$ws.Range("F1").ColumnWidth = 19.5
$ws.Range("G1").ColumnWidth = 17.666666666666668

# --- Turn the range into an Excel Table ("Table1") --------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:G7"), 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleMedium2"
$lo.ShowAutoFilterDropDown = $false

# --- Make Sheet1 the active sheet / selection, matching the saved state ----
$ws.Activate()
$ws.Range("A1:G7").Select()
